$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), matching the formatting of the
# existing header cells (e.g. G1) which use style index 1
# (bold font, border, centered alignment).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data value for the "Save" column in row 2.
$ws.Range("H2").Value = 1
